$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 37491
$ws.Range("D2").Value = 54233007
$ws.Range("C3").Value = 90518
$ws.Range("D3").Value = 132708406
$ws.Range("C4").Value = 31027
$ws.Range("D4").Value = 45951955
$ws.Range("C5").Value = 8652
$ws.Range("D5").Value = 12861235
$ws.Range("C6").Value = 1978
$ws.Range("D6").Value = 2939506
$ws.Range("C11").Value = 41013
$ws.Range("D11").Value = 55667487
$ws.Range("C12").Value = 9603
$ws.Range("D12").Value = 13890569
$ws.Range("C13").Value = 25831
$ws.Range("D13").Value = 37884302
$ws.Range("C14").Value = 8288
$ws.Range("D14").Value = 12301263
$ws.Range("C19").Value = 10168
$ws.Range("D19").Value = 13470847
$ws.Range("C20").Value = 13320
$ws.Range("D20").Value = 19237976
$ws.Range("C21").Value = 31529
$ws.Range("D21").Value = 46275398
$ws.Range("C22").Value = 10196
$ws.Range("D22").Value = 15157055
$ws.Range("C23").Value = 2620
$ws.Range("D23").Value = 3896155
$ws.Range("C24").Value = 503
$ws.Range("D24").Value = 748592
$ws.Range("C26").Value = 11615
$ws.Range("D26").Value = 15521381
$ws.Range("C27").Value = 7607
$ws.Range("D27").Value = 11022003
$ws.Range("C28").Value = 22375
$ws.Range("D28").Value = 32842554
$ws.Range("C29").Value = 7772
$ws.Range("D29").Value = 11566302
$ws.Range("C30").Value = 1952
$ws.Range("D30").Value = 2912499
$ws.Range("C31").Value = 365
$ws.Range("D31").Value = 544915
$ws.Range("C33").Value = 8254
$ws.Range("D33").Value = 10908103
$ws.Range("C34").Value = 3206
$ws.Range("D34").Value = 4627280
$ws.Range("C35").Value = 7759
$ws.Range("D35").Value = 11332446
$ws.Range("C36").Value = 3161
$ws.Range("D36").Value = 4684461
$ws.Range("C37").Value = 821
$ws.Range("D37").Value = 1222823
$ws.Range("C40").Value = 2443
$ws.Range("D40").Value = 3302027
$ws.Range("C41").Value = 17127
$ws.Range("D41").Value = 24771633
$ws.Range("C42").Value = 50835
$ws.Range("D42").Value = 74534046
$ws.Range("C43").Value = 18927
$ws.Range("D43").Value = 28115688
$ws.Range("C44").Value = 5581
$ws.Range("D44").Value = 8311478
$ws.Range("C45").Value = 1191
$ws.Range("D45").Value = 1777045
$ws.Range("C49").Value = 16584
$ws.Range("D49").Value = 22092496
$ws.Range("C50").Value = 1985
$ws.Range("D50").Value = 2879982
$ws.Range("C51").Value = 6790
$ws.Range("D51").Value = 9983833
$ws.Range("C52").Value = 2324
$ws.Range("D52").Value = 3470918
$ws.Range("C53").Value = 748
$ws.Range("D53").Value = 1117305
$ws.Range("C54").Value = 182
$ws.Range("D54").Value = 269833
$ws.Range("C56").Value = 6727
$ws.Range("D56").Value = 9270974
$ws.Range("C57").Value = 918
$ws.Range("D57").Value = 1347084
$ws.Range("C58").Value = 2290
$ws.Range("D58").Value = 3394817
$ws.Range("C59").Value = 911
$ws.Range("D59").Value = 1356001
$ws.Range("C60").Value = 315
$ws.Range("D60").Value = 472258
$ws.Range("C63").Value = 1343
$ws.Range("D63").Value = 1892985
$ws.Range("C64").Value = 15255
$ws.Range("D64").Value = 22039010
$ws.Range("C65").Value = 44467
$ws.Range("D65").Value = 65076781
$ws.Range("C66").Value = 15631
$ws.Range("D66").Value = 23233714
$ws.Range("C67").Value = 4550
$ws.Range("D67").Value = 6777292
$ws.Range("C68").Value = 913
$ws.Range("D68").Value = 1357668
$ws.Range("C72").Value = 15020
$ws.Range("D72").Value = 19810924
$ws.Range("C73").Value = 50849
$ws.Range("D73").Value = 74000345
$ws.Range("C74").Value = 144761
$ws.Range("D74").Value = 213284497
$ws.Range("C75").Value = 63156
$ws.Range("D75").Value = 94111434
$ws.Range("C76").Value = 20167
$ws.Range("D76").Value = 30131317
$ws.Range("C77").Value = 4763
$ws.Range("D77").Value = 7116223
$ws.Range("C84").Value = 50345
$ws.Range("D84").Value = 68531331
$ws.Range("C85").Value = 4548
$ws.Range("D85").Value = 6588941
$ws.Range("C86").Value = 11474
$ws.Range("D86").Value = 16858056
$ws.Range("C87").Value = 3857
$ws.Range("D87").Value = 5748406
$ws.Range("C88").Value = 1336
$ws.Range("D88").Value = 1996489
$ws.Range("C89").Value = 285
$ws.Range("D89").Value = 425012
$ws.Range("C92").Value = 5334
$ws.Range("D92").Value = 7170465
$ws.Range("C93").Value = 1572
$ws.Range("D93").Value = 2262932
$ws.Range("C94").Value = 5099
$ws.Range("D94").Value = 7512029
$ws.Range("C95").Value = 1928
$ws.Range("D95").Value = 2872446
$ws.Range("C96").Value = 684
$ws.Range("D96").Value = 1024960
$ws.Range("C97").Value = 179
$ws.Range("D97").Value = 267613
$ws.Range("C100").Value = 3489
$ws.Range("D100").Value = 4623032
$ws.Range("C101").Value = 590
$ws.Range("D101").Value = 878664
$ws.Range("C102").Value = 344
$ws.Range("D102").Value = 513530
$ws.Range("C103").Value = 129
$ws.Range("D103").Value = 193500
$ws.Range("C106").Value = 10700
$ws.Range("D106").Value = 15529772
$ws.Range("C107").Value = 29066
$ws.Range("D107").Value = 42710014
$ws.Range("C108").Value = 9741
$ws.Range("D108").Value = 14486213
$ws.Range("C110").Value = 485
$ws.Range("D110").Value = 722546
$ws.Range("C113").Value = 9725
$ws.Range("D113").Value = 12855726
$ws.Range("C114").Value = 30212
$ws.Range("D114").Value = 43574570
$ws.Range("C115").Value = 65811
$ws.Range("D115").Value = 96322794
$ws.Range("C116").Value = 21270
$ws.Range("D116").Value = 31609880
$ws.Range("C117").Value = 6028
$ws.Range("D117").Value = 8981826
$ws.Range("C123").Value = 25678
$ws.Range("D123").Value = 34312660
$ws.Range("C124").Value = 35709
$ws.Range("D124").Value = 51545193
$ws.Range("C125").Value = 76312
$ws.Range("D125").Value = 111606054
$ws.Range("C126").Value = 23716
$ws.Range("D126").Value = 35201035
$ws.Range("C127").Value = 6346
$ws.Range("D127").Value = 9430051
$ws.Range("C128").Value = 1220
$ws.Range("D128").Value = 1814411
$ws.Range("C132").Value = 31522
$ws.Range("D132").Value = 41881184
$ws.Range("C133").Value = 13155
$ws.Range("D133").Value = 19044258
$ws.Range("C134").Value = 32156
$ws.Range("D134").Value = 47233629
$ws.Range("C135").Value = 11431
$ws.Range("D135").Value = 16984552
$ws.Range("C136").Value = 2944
$ws.Range("D136").Value = 4388805
$ws.Range("C137").Value = 492
$ws.Range("D137").Value = 731990
$ws.Range("C140").Value = 10759
$ws.Range("D140").Value = 14354603
$ws.Range("C141").Value = 34791
$ws.Range("D141").Value = 50251654
$ws.Range("C142").Value = 80746
$ws.Range("D142").Value = 118309013
$ws.Range("C143").Value = 24236
$ws.Range("D143").Value = 36012622
$ws.Range("C144").Value = 6353
$ws.Range("D144").Value = 9480208
$ws.Range("C145").Value = 1420
$ws.Range("D145").Value = 2112230
$ws.Range("C146").Value = 79
$ws.Range("D146").Value = 118130
$ws.Range("C148").Value = 28962
$ws.Range("D148").Value = 39100111
